$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the symptom-group label in A8
$ws.Range("A8").Value = "Symptom – Musculoskeletal"

# Updated data values (rows 2-12, columns B-H)
$data = @{
    2  = @(5.3, 4.3, 3.4, 7.5, 6.2, 3.9, 4.1)
    3  = @(17.7, 12.6, 9, 18, 10.5, 5.5, 10.9)
    4  = @(8.300000000000001, 5.5, 8, 6.4, 7.2, 7.1, 6.8)
    5  = @(9.4, 7.2, 8.1, 11.3, 8, 9, 9.1)
    6  = @(15.4, 10.6, 9.800000000000001, 12.8, 14.9, 10.6, 12.2)
    7  = @(4.1, 4.3, 4.5, 4.1, 6.2, 6.7, 5.3)
    8  = @(2.6, 2, 2.2, 2.9, 1.4, 1.2, 1.8)
    9  = @(5.6, 11.3, 11.2, 9.800000000000001, 13.4, 13.7, 12.7)
    10 = @(22.6, 38.2, 40, 21.1, 27.9, 40, 34.4)
    11 = @(3.8, 2, 2.2, 3.3, 3.3, 1.2, 1.9)
    12 = @(5.3, 2.1, 1.6, 2.9, 1.1, 1.2, 0.9)
}

$cols = @("B", "C", "D", "E", "F", "G", "H")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cellRef = "$($cols[$i])$row"
        $ws.Range($cellRef).Value = $values[$i]
    }
}
